$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '244.16'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-0.53%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '26.53'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '3.89%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.133'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '0.05%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05608'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.43%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.468'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-0.27%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8193'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '0.18%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8333'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-0.88%'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1326'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-0.69%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06927'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-0.49%'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.02894'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '1.10%'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09380'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.04%'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001510'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.55%'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'One'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.01002'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '1,578.86%'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.006145'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-0.08%'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.650'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '3.35%'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.020'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-0.09%'
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.301'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '13.78%'
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3112'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-2.11%'
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.03081'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-4.25%'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-1.47%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.752'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '0.10%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04584'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-2.05%'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-2.42%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001228'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-1.45%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004494'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.00009597'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-1.04%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03640'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-0.52%'
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006168'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '0.48%'
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1051'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-0.19%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002399'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-5.12%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008116'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '4.41%'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '0.66%'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.00%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002503'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '17.86%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002099'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '0.00%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0001999'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.00%'
